{"js": "// This script updates the \"Sprint No.\" and \"Review Date\" values in the\n// document's header metadata table:\n//   Sprint No.:   1         -> 2\n//   Review Date:  02/09/18  -> 02/21/18\n//\n// The values are looked up precisely via the table/row/cell structure so\n// that similarly-looking text elsewhere in the document (e.g. the \"1.\"\n// in the \"1.  License\" heading) is left untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row 1 (\"Reviewer's Name\" / \"Sprint No.\") -> last cell holds the sprint number.\nconst sprintRow = rows.items[1];\nconst sprintCells = sprintRow.cells;\nsprintCells.load(\"items\");\nawait context.sync();\n\nconst sprintCell = sprintCells.items[sprintCells.items.length - 1];\nconst sprintResults = sprintCell.body.search(\"1\", { matchCase: true, matchWholeWord: true });\nsprintResults.load(\"items/text\");\nawait context.sync();\n\nfor (const r of sprintResults.items) {\n  if (r.text === \"1\") {\n    r.insertText(\"2\", \"Replace\");\n  }\n}\nawait context.sync();\n\n// Row 2 (\"Review Date\") -> second cell holds the date value.\nconst dateRow = rows.items[2];\nconst dateCells = dateRow.cells;\ndateCells.load(\"items\");\nawait context.sync();\n\nconst dateCell = dateCells.items[1];\nconst dateResults = dateCell.body.search(\"02/09/18\", { matchCase: true });\ndateResults.load(\"items/text\");\nawait context.sync();\n\nfor (const r of dateResults.items) {\n  r.insertText(\"02/21/18\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Updates the \"Sprint No.\" and \"Review Date\" values in the document's\n# header metadata table:\n#   Sprint No.:   1         -> 2\n#   Review Date:  02/09/18  -> 02/21/18\n#\n# The values are located precisely via the table/row/cell structure so that\n# similarly-looking text elsewhere in the document (e.g. the \"1.\" in the\n# \"1.  License\" heading) is left untouched. Assigning .Text directly (rather\n# than using Find/Replace) keeps the existing run/formatting intact.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Row 2 (\"Reviewer's Name\" / \"Sprint No.\") -> last cell holds the sprint number.\n$sprintCell = $table.Cell(2, $table.Rows.Item(2).Cells.Count)\n$sprintRange = $sprintCell.Range\n$sprintRange.MoveEnd(1, -1) | Out-Null  # exclude the trailing cell-mark character\n$sprintRange.Text = \"2\"\n\n# Row 3 (\"Review Date\") -> second cell holds the date value.\n$dateCell = $table.Cell(3, 2)\n$dateRange = $dateCell.Range\n$dateRange.MoveEnd(1, -1) | Out-Null  # exclude the trailing cell-mark character\n$dateRange.Text = \"02/21/18\"\n"}
